# Adds "requiredMsg" / "invalidMsg" validation-message columns (C, D) to the
# inValidLoginData sheet: a header in row 1, plus "Required" / "Invalid
# credentials" alongside every existing data row (including a brand-new,
# otherwise-empty row 6 between the current rows 5 and 7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C: header, then "Required" for every data row (2-9). Row 6 previously
# had no data at all; it now gets C/D only.
$ws.Cells.Item(1, 3).Value = "requiredMsg"
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value = "Required"
}

# Column D: data rows first, header last (matches the authored edit order).
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 4).Value = "Invalid credentials"
}
$ws.Cells.Item(1, 4).Value = "invalidMsg"

# New column widths for the validation-message columns (closest values the
# Excel object model's pixel-quantized ColumnWidth can reach to 17.140625 /
# 18.42578125 respectively)
$ws.Columns.Item(3).ColumnWidth = 16.3
$ws.Columns.Item(4).ColumnWidth = 17.6

# Match the saved selection: D2:D9 with D2 active
$ws.Range("D2:D9").Select()
